$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Copy style/formatting from column CB into new column CC for all used rows
$ws.Range("CB1:CB25").Copy($ws.Range("CC1:CC25"))

# Header
$ws.Range("CC1").Value = "02-sep"

# Data rows (02-sep prices)
$ws.Range("CC2").Value = 33.61
$ws.Range("CC3").Value = 36.79
$ws.Range("CC4").Value = 52.57
$ws.Range("CC5").Value = 40.16
$ws.Range("CC6").Value = 42.55
$ws.Range("CC7").Value = 12.47
$ws.Range("CC8").Value = 33.44
$ws.Range("CC9").Value = 34.28
$ws.Range("CC10").Value = 34.8
$ws.Range("CC11").Value = 39
$ws.Range("CC12").Value = 7.96
$ws.Range("CC13").Value = 1.2
$ws.Range("CC14").Value = 1.5
$ws.Range("CC15").Value = 0.65
$ws.Range("CC16").Value = 0.21
$ws.Range("CC17").Value = 1.5
$ws.Range("CC18").Value = 8.210000000000001
$ws.Range("CC19").Value = 14
$ws.Range("CC20").Value = 32.89
$ws.Range("CC21").Value = 84
$ws.Range("CC22").Value = 53.46
$ws.Range("CC23").Value = 34.98
$ws.Range("CC24").Value = 42.89
$ws.Range("CC25").Value = 31.06
